# "first results android emulator"
# Rename the first sheet ("Feuil1") to "Real" (it now holds the first batch
# of real-device results, as opposed to the Genymotion emulator sheet),
# make sure it's the active sheet/tab, and move the selection to where the
# author last left off (L18) while working on it.

$wb = $excel.ActiveWorkbook

$wsReal = $wb.Worksheets.Item(1)
$wsReal.Name = "Real"

$wsReal.Activate()
$wsReal.Range("L18").Select()
